# Load-test fixture data: expand the single-month sample (rows 3-8) into a full
# per-fuel-type activity table (rows 3-18) spanning multiple periods, matching the
# "Actividad / Tipo de Consumo / Consumo / Unidad / periodicidad / periodoDeImputacion"
# columns already defined in the header (rows 1-2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting first, while the template's original per-row formats (plain vs
# wrapped text in columns A/B) are still in place on rows 3-8, so we can clone
# them onto the rows whose content is changing/being added. B4 is the template's
# last surviving 'wrapped' cell once the row-3 'plain' format below overwrites the
# others, so copy the wrapped look first.
$ws.Range("B4").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 9 & 10 are brand-new (the template only pre-formats rows 11+ with a 15.75pt
# row height); give them the same 15pt height the rest of the 'Combustion Fija'
# block already had in the template file's row-height runs.
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15

# --- Cell values -----------------------------------------------------------
$ws.Cells.Item(3,1).Value = "Combustión Fija"
$ws.Cells.Item(3,2).Value = "Gas Natural"
$ws.Cells.Item(3,3).Value = 200
$ws.Cells.Item(3,4).Value = "m3"
$ws.Cells.Item(3,5).Value = "Mensual"
$ws.Cells.Item(3,6).Value = 44531

$ws.Cells.Item(4,1).Value = "Combustión Fija"
$ws.Cells.Item(4,2).Value = "Diesel/Gasoil"
$ws.Cells.Item(4,3).Value = 500
$ws.Cells.Item(4,4).Value = "lt"
$ws.Cells.Item(4,5).Value = "Anual"
$ws.Cells.Item(4,6).Value = 44562

$ws.Cells.Item(5,1).Value = "Combustión Fija"
$ws.Cells.Item(5,2).Value = "Kerosene"
$ws.Cells.Item(5,3).Value = 100
$ws.Cells.Item(5,4).Value = "lt"
$ws.Cells.Item(5,5).Value = "Mensual"
$ws.Cells.Item(5,6).Value = 44652

$ws.Cells.Item(6,1).Value = "Combustión Fija"
$ws.Cells.Item(6,2).Value = "Fuel Oil"
$ws.Cells.Item(6,3).Value = 750
$ws.Cells.Item(6,4).Value = "lt"
$ws.Cells.Item(6,5).Value = "Anual"
$ws.Cells.Item(6,6).Value = 44287

$ws.Cells.Item(7,1).Value = "Combustión Fija"
$ws.Cells.Item(7,2).Value = "Nafta"
$ws.Cells.Item(7,3).Value = 1800
$ws.Cells.Item(7,4).Value = "lt"
$ws.Cells.Item(7,5).Value = "Mensual"
$ws.Cells.Item(7,6).Value = 44317

$ws.Cells.Item(8,1).Value = "Combustión Fija"
$ws.Cells.Item(8,2).Value = "Carbón"
$ws.Cells.Item(8,3).Value = 245
$ws.Cells.Item(8,4).Value = "kg"
$ws.Cells.Item(8,5).Value = "Mensual"
$ws.Cells.Item(8,6).Value = 44348

$ws.Cells.Item(9,1).Value = "Combustión Fija"
$ws.Cells.Item(9,2).Value = "Carbón de leña"
$ws.Cells.Item(9,3).Value = 650
$ws.Cells.Item(9,4).Value = "lt"
$ws.Cells.Item(9,5).Value = "Mensual"
$ws.Cells.Item(9,6).Value = 44378

$ws.Cells.Item(10,1).Value = "Combustión Fija"
$ws.Cells.Item(10,2).Value = "Leña"
$ws.Cells.Item(10,3).Value = 850
$ws.Cells.Item(10,4).Value = "lt"
$ws.Cells.Item(10,5).Value = "Anual"
$ws.Cells.Item(10,6).Value = 44409

$ws.Cells.Item(11,1).Value = "Combustión Móvil"
$ws.Cells.Item(11,2).Value = "Combustible Consumido - Gasoil"
$ws.Cells.Item(11,3).Value = 200
$ws.Cells.Item(11,4).Value = "lt"
$ws.Cells.Item(11,5).Value = "Mensual"
$ws.Cells.Item(11,6).Value = 44743

$ws.Cells.Item(12,1).Value = "Combustión Móvil"
$ws.Cells.Item(12,2).Value = "Combustible Consumido - GNC"
$ws.Cells.Item(12,3).Value = 1200
$ws.Cells.Item(12,4).Value = "lt"
$ws.Cells.Item(12,5).Value = "Mensual"
$ws.Cells.Item(12,6).Value = 44470

$ws.Cells.Item(13,1).Value = "Combustión Móvil"
$ws.Cells.Item(13,2).Value = "Combustible Consumido - Nafta"
$ws.Cells.Item(13,3).Value = 1400
$ws.Cells.Item(13,4).Value = "lt"
$ws.Cells.Item(13,5).Value = "Mensual"
$ws.Cells.Item(13,6).Value = 44501

$ws.Cells.Item(14,1).Value = "Electricidad Adquirida y Consumida"
$ws.Cells.Item(14,2).Value = "Electricidad"
$ws.Cells.Item(14,3).Value = 700
$ws.Cells.Item(14,4).Value = "kwh"
$ws.Cells.Item(14,5).Value = "Anual"
$ws.Cells.Item(14,6).Value = 44287

$ws.Cells.Item(15,1).Value = "Logística de productos y residuos"
$ws.Cells.Item(15,2).Value = "Categoría de producto transportado"
$ws.Cells.Item(15,3).Value = 100
$ws.Cells.Item(15,4).Value = "-"
$ws.Cells.Item(15,5).Value = "Mensual"
$ws.Cells.Item(15,6).Value = 44682

$ws.Cells.Item(16,1).Value = "Logística de productos y residuos"
$ws.Cells.Item(16,2).Value = "Medio de Transporte"
$ws.Cells.Item(16,3).Value = 650
$ws.Cells.Item(16,4).Value = "-"
$ws.Cells.Item(16,5).Value = "Mensual"
$ws.Cells.Item(16,6).Value = 44593

$ws.Cells.Item(17,1).Value = "Logística de productos y residuos"
$ws.Cells.Item(17,2).Value = "Distancia Medio Recorrida"
$ws.Cells.Item(17,3).Value = 750
$ws.Cells.Item(17,4).Value = "km"
$ws.Cells.Item(17,5).Value = "Anual"
$ws.Cells.Item(17,6).Value = 44621

$ws.Cells.Item(18,1).Value = "Logística de productos y residuos"
$ws.Cells.Item(18,2).Value = "Peso Total Transportado"
$ws.Cells.Item(18,3).Value = 980
$ws.Cells.Item(18,4).Value = "kg"
$ws.Cells.Item(18,5).Value = "Mensual"
$ws.Cells.Item(18,6).Value = 44713

[void]$ws.Range("F3:F18").Select()
